$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 10488.5
$ws.Range("B7").Value = 10408.36
$ws.Range("C7").Value = 107.96
$ws.Range("D7").Value = 108.79
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 0.77
$ws.Range("G7").Value = 42609.48847222222
$ws.Range("G7").NumberFormat = "m/d/yy h:mm"
$ws.Range("H7").Value = $true
